$d = $word.ActiveDocument

# --- Remove the "[post around Jan 3/4th so gets 3D printing article as homepage]" note ---
# Clear the text content of paragraph 4 (leaving it as an empty paragraph, like the
# surrounding placeholder paragraphs) rather than deleting the paragraph mark itself.
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4content = $d.Range($r4.Start, $r4.End - 1)
$r4content.Delete()

# The two following blank paragraphs are no longer needed now that paragraph 4 itself
# is blank (the "Each day you go on..." paragraph should immediately follow two blanks).
$d.Paragraphs.Item(5).Range.Delete()
$d.Paragraphs.Item(5).Range.Delete()

# --- Extend the "A small LLM summarizes..." paragraph and split off the "[2/5]" tag ---
$p9 = $d.Paragraphs.Item(9)
$r9 = $p9.Range
$r9content = $d.Range($r9.Start, $r9.End - 1)
$r9content.Delete()

$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertBefore("A small LLM summarizes what was written back then, and how the future unfolded to see how predictions matched reality. So you can look back at scientific ideas that were once thought to be potential breakthroughs and see what real-world impact they actually had. ")

# New paragraph holding the "[3/5]" counter, replacing the old trailing "[2/5]" run.
$p9 = $d.Paragraphs.Item(9)
$p9.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertBefore("[3/5]")

# --- Add the new closing paragraphs (still before the document's final blank paragraph) ---
$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs.Item(12)
$p12.Range.InsertParagraphAfter()
$p13 = $d.Paragraphs.Item(13)
$p13.Range.InsertBefore("Many years of articles are already up on the website, and more are being added regularly, so you can jump across eras all the way from the completion of the Human Genome Project to the first days of the COVID pandemic.  [4/5]")

$p13 = $d.Paragraphs.Item(13)
$p13.Range.InsertParagraphAfter()
$p14 = $d.Paragraphs.Item(14)
$p14.Range.InsertParagraphAfter()
$p15 = $d.Paragraphs.Item(15)
$p15.Range.InsertBefore("Be warned that the LLM is not perfect in factual accuracy, especially on very niche pages. ")

$p15 = $d.Paragraphs.Item(15)
$p15.Range.InsertParagraphAfter()
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertBefore("Overall it’s done okay so far considering it’s a small side-project. If there is conitnued interest in this, I’ll keep working to improve the quality & accuracy of the analysis model! [5/5]")
